# Update "想去人数" (F column) figures for two sheets: 展览 and 全部类型
$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 115
    4  = 1623
    8  = 11516
    15 = 12382
    16 = 13085
    24 = 121
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
